# ModificarDescipcionFacturaInscrita.xlsx update
# - updates the "fechaFactura" data value on the Datos sheet (N2)
# - moves the active selection on the Datos sheet to N7 (as recorded by the
#   resave that produced this commit)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Update the stored date-serial value in N2 (fechaFactura column)
$ws.Range("N2").Value = 65468

# Leave the sheet with the same active-cell selection captured in the diff
[void]$ws.Range("N7").Select()
